# Roboflow Annotation Report 7/22/2025
# Append a new weekly-progress row to the Table1 listobject on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Table1")

# Adding a ListRow grows the table (ref + autoFilter) and the sheet
# dimension automatically, right below the current last data row.
$newListRow = $tbl.ListRows.Add()
$rowIndex = $newListRow.Range.Row
$srcRow = $rowIndex - 1

# Match the look of the row above (borders/number format/font) by copying
# its formatting into the freshly-added row.
$ws.Range("D" + $srcRow + ":J" + $srcRow).Copy()
$ws.Range("D" + $rowIndex + ":J" + $rowIndex).PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item($rowIndex).RowHeight = $ws.Rows.Item($srcRow).RowHeight
$excel.CutCopyMode = 0

$ws.Cells.Item($rowIndex, 4).Value = "22/7/2026"
$ws.Cells.Item($rowIndex, 5).Value = 396
$ws.Cells.Item($rowIndex, 6).Value = 934
$ws.Cells.Item($rowIndex, 7).Value = 0
$ws.Cells.Item($rowIndex, 8).Value = 0
$ws.Cells.Item($rowIndex, 9).Value = 1012
$ws.Cells.Item($rowIndex, 10).Value = "N/A"

$ws.Range("F67").Select()
